$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook / window view tweaks -----------------------------------------
# Tab ratio (split between sheet tabs and horizontal scrollbar) change.
try {
    $wb.Windows.Item(1).TabRatio = 0.076
} catch {
}

# Scroll the view back to the top-left and move the selection.
try {
    $excel.ActiveWindow.ScrollRow = 1
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("S25").Select() | Out-Null

# --- Newly synthesized "alu" Yosys-ABC results (row 15 & 16, cols D:I) -----
$ws.Range("D15").Value = 41
$ws.Range("E15").Value = 111
$ws.Range("F15").Formula = "=D15+E15"
$ws.Range("G15").Value = 195
$ws.Range("H15").Formula = "=D15+E15+G15"
$ws.Range("I15").Formula = "=F15+(G15*5)"

$ws.Range("D16").Value = 211
$ws.Range("E16").Value = 346
$ws.Range("F16").Formula = "=D16+E16"
$ws.Range("G16").Value = 685
$ws.Range("H16").Formula = "=D16+E16+G16"
$ws.Range("I16").Formula = "=F16+(G16*5)"

# --- Apply the "#,##0" thousands-separator format to the whole Synopsys ----
# block (columns P:U, rows 5:16). This also materialises the previously
# blank separator-row cells (rows 7, 10, 11, 14) with that same format.
$ws.Range("P5:U16").NumberFormat = "#,##0"

# The new alu formula columns (F, H, I) get the same numeric format as the
# rest of the "Total"/"Weight" columns elsewhere in the sheet.
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("H15").NumberFormat = "#,##0"
$ws.Range("I15").NumberFormat = "#,##0"
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("H16").NumberFormat = "#,##0"
$ws.Range("I16").NumberFormat = "#,##0"

# Row 14 also picks up the same "#,##0" styling on F/H/I (mirrors the other
# blank separator rows for the Yosys-ABC "Total"/"Weight" columns).
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("H14").NumberFormat = "#,##0"
$ws.Range("I14").NumberFormat = "#,##0"

# --- Clear the inherited bold column formatting on a few blank cells so ---
# they match the plain "Normal" style used throughout the rest of column A/B.
$ws.Range("A13").Font.Bold = $false
$ws.Range("A14").Font.Bold = $false
$ws.Range("B14").Font.Bold = $false
